$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Israel Premier League")

# Rows 164 and 165 had their match-data columns (B, E:AD) swapped between them.
$row164 = @{
    "B" = 6799984
    "E" = "Hapoel Bnei Sakhnin"
    "F" = "MS Ashdod"
    "H" = 1
    "J" = 1
    "K" = "D"
    "L" = 2
    "M" = 3.25
    "N" = 3.25
    "O" = 2
    "P" = 3.25
    "Q" = 3.25
    "S" = 1.825
    "T" = 2.025
    "U" = 2.25
    "V" = 1.925
    "W" = 1.925
    "X" = -1
    "Y" = 2.25
    "AA" = -0.5
    "AB" = 0.5125
    "AC" = -0.5
    "AD" = 0.4625
}

$row165 = @{
    "B" = 6799986
    "E" = "Hapoel Jerusalem FC"
    "F" = "Maccabi Bnei Raina"
    "H" = 0
    "J" = 0
    "K" = "H"
    "L" = 2.2
    "M" = 3.2
    "N" = 2.9
    "O" = 2.3
    "P" = 3.1
    "Q" = 2.8
    "S" = 2.1
    "T" = 1.775
    "U" = 2
    "V" = 2.025
    "W" = 1.825
    "X" = 1.3
    "Y" = -1
    "AA" = 1.1
    "AB" = -1
    "AC" = -1
    "AD" = 0.825
}

# Rows 211 and 212 had their match-data columns (B, E:AD) swapped between them.
$row211 = @{
    "B" = 8016170
    "E" = "Hapoel TelAviv"
    "F" = "Beitar Jerusalem"
    "H" = 5
    "I" = 0
    "J" = 1
    "K" = "A"
    "L" = 2.875
    "M" = 3.2
    "N" = 2.5
    "O" = 3.2
    "P" = 3.2
    "Q" = 2.3
    "R" = 0.25
    "S" = 1.875
    "T" = 1.975
    "U" = 2.5
    "V" = 1.975
    "W" = 1.875
    "X" = -1
    "Y" = -1
    "Z" = 1.3
    "AA" = -1
    "AB" = 0.9750000000000001
    "AC" = 0.9750000000000001
    "AD" = -1
}

$row212 = @{
    "B" = 8016155
    "E" = "MS Ashdod"
    "F" = "Hapoel Petah Tikva"
    "H" = 0
    "I" = 1
    "J" = 0
    "K" = "H"
    "L" = 2.1
    "M" = 3.3
    "N" = 3.5
    "O" = 2.2
    "P" = 3.25
    "Q" = 3.3
    "R" = -0.25
    "S" = 1.9
    "T" = 1.95
    "U" = 2.25
    "V" = 1.925
    "W" = 1.925
    "X" = 1.2
    "Y" = -1
    "Z" = -1
    "AA" = 0.8999999999999999
    "AB" = -1
    "AC" = -1
    "AD" = 0.925
}

function Set-RowValues($ws, [int]$rowNum, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$rowNum").Value = $values[$col]
    }
}

Set-RowValues $ws 164 $row164
Set-RowValues $ws 165 $row165
Set-RowValues $ws 211 $row211
Set-RowValues $ws 212 $row212
